$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row at 10 (BirthdayList_Index_ResourceName), shifting everything else down by one.
$ws.Rows.Item(10).Insert()

# New row 10: BirthdayList_Index_ResourceName
$ws.Range("A10").Value = "BirthdayList_Index_ResourceName"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 'Column index of "Resource Name" in BirthdayList datatable'

# Row 17 (was UpcomingBirthday_IncludeTeam) - update team name value
$ws.Range("B17").Value = "Management I"

# Row 20 (Email_Subject) - include celebrant name placeholder in subject, update description
$ws.Range("B20").Value = "Testing - upcoming birthday for: {0}"
$ws.Range("C20").Value = "Email subject to be sent to the team members and management. Where 0 is the value for the celebrant's name."

# Row 21 (Email_Body) - add hyperlink to sharepoint list
$ws.Range("B21").Value = 'Testing - Please send your birthday greetings <a href="https://ts.accenture.com/sites/Chevron_PDC/CET/Lists/CET%20Birthday%20Greetings/AllItems.aspx">here</href>'

# Update view: scroll + selection
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B22").Select()
